# LOM3049.xlsx edit script
# Reassigns the "Docentes responsaveis" / "Programa" / assessment rows so the
# sheet ends up with one fewer content row (25 -> 24), matching the target
# shared-string/content layout described in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value changes -----------------------------------------------

# Row 10 (Objetivos:): long paragraph replaced by the professor name.
$ws.Range("B10").Value = "5840521 - Rosa Ana Conte"
$ws.Range("C10").Value = "5840521 - Rosa Ana Conte"

# Row 13 used to hold only the professor name in B/C with blank A.
# It becomes the "Programa resumido:" row, now showing "Semestral".
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 14 becomes "Short syllabus:" with no B/C value.
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()

# Row 15 becomes "Programa:" and picks up the "01/01/2021" text in B/C.
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"

# Row 16 becomes "Syllabus:" with no B/C value.
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()

# Row 17 becomes "Avaliação:".
$ws.Range("A17").Value = "Avaliação:"

# Row 18 becomes "Método:" and picks up the professor-name text in B/C.
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840521 - Rosa Ana Conte"
$ws.Range("C18").Value = "5840521 - Rosa Ana Conte"

# Row 19 becomes "Critério:" (text in B/C unchanged - the "Serão
# realizadas 2 avaliações..." paragraph).
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Serão realizadas 2 avaliações, com questões abrangendo problemas práticos e conceituais. A 1a. avaliação terá peso 1 e a 2a. avaliação terá peso 2. A nota será a média ponderada das 2 avaliações."
$ws.Range("C19").Value = "Serão realizadas 2 avaliações, com questões abrangendo problemas práticos e conceituais. A 1a. avaliação terá peso 1 e a 2a. avaliação terá peso 2. A nota será a média ponderada das 2 avaliações."

# Row 20 becomes "Norma de recuperação:" (text in B/C unchanged - the
# "Serão aplicadas duas avaliações escritas..." paragraph).
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Serão aplicadas duas avaliações escritas (P1, com peso 1 e P2, com peso 2) que comporão a nota final (NF). A nota final será calculada através da expressão: NF = (P1 + P2)/3."
$ws.Range("C20").Value = "Serão aplicadas duas avaliações escritas (P1, com peso 1 e P2, com peso 2) que comporão a nota final (NF). A nota final será calculada através da expressão: NF = (P1 + P2)/3."

# Row 21 becomes "Bibliografia:" (text in B/C unchanged - the
# "Para a recuperação será realizada uma prova..." paragraph).
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2."
$ws.Range("C21").Value = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2."

# Row 22 becomes "Requisitos:" with no B/C value (old bibliography
# paragraph is dropped entirely).
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()

# Row 23 loses its "Requisitos:" label (moved to row 22) and now holds
# the first prerequisite line in B/C.
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOB1004 -  Cálculo II  (Requisito)`n"
$ws.Range("C23").Value = "LOB1004 -  Cálculo II  (Requisito)`n"

# Row 24 now holds the second prerequisite line in B/C.
$ws.Range("B24").Value = "LOB1019 -  Física II  (Requisito)`n"
$ws.Range("C24").Value = "LOB1019 -  Física II  (Requisito)`n"

# Row 25 no longer exists - remove it and shift everything below (none) up.
$ws.Rows("25").Delete()

# --- Row height changes -------------------------------------------------

$ws.Rows("11").AutoFit()
$ws.Rows("13").RowHeight = 60
$ws.Rows("14").AutoFit()
$ws.Rows("15").RowHeight = 120
$ws.Rows("16").AutoFit()
$ws.Rows("17").AutoFit()
$ws.Rows("18").RowHeight = 60
$ws.Rows("21").RowHeight = 120
$ws.Rows("22").AutoFit()
$ws.Rows("23").RowHeight = 30
